$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Update the note text in U8 (append a trailing period + space to the sentence)
$ws.Range("U8").Value = "Para el periodo reportado no se han financiado Estudios con recursos públicos. "

# Update the reporting period dates (row 8): Q1 2022 -> Q2 2022
$ws.Range("B8").Value = (Get-Date -Year 2022 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C8").Value = (Get-Date -Year 2022 -Month 6 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("S8").Value = (Get-Date -Year 2022 -Month 7 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("T8").Value = (Get-Date -Year 2022 -Month 7 -Day 11 -Hour 0 -Minute 0 -Second 0)

# Row 3: grow the row to fit the wrapped description text in G3:I3
$ws.Rows.Item(3).RowHeight = 35.25

# Strip the full box border from G3, then restore only the left edge, and
# turn wrap-text on for the merged description cell (G3:I3)
$ws.Range("G3").Borders.LineStyle = -4142
$ws.Range("G3").Borders.Item(7).LineStyle = 1
$ws.Range("G3").WrapText = $true
$ws.Range("H3").WrapText = $true
$ws.Range("I3").WrapText = $true

# Update the active selection to match the saved view state
$ws.Range("U15").Select()
